$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update status ("estado") values to "terminado" for the affected tasks
$ws.Range("B50").Value = "terminado"
$ws.Range("B52").Value = "terminado"
$ws.Range("B55").Value = "terminado"

# Update the visible scroll position / active cell selection to match the new view
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("B54").Select()
